# HRMS Regression test script First Commit
#
# Applies to EmployeeAttendanceCalender.xlsx:
#  1. Adds a new (blank) trailing worksheet "Sheet5".
#  2. Renames the employee on the "AttendanceApproval" sheet (B2) from
#     "Ramnaresh  Sarwan" to "F Dinesh", and gives that cell an explicit
#     black font colour (new font/style, matching the new <fonts>/<cellXfs>
#     entries in the target workbook).
#  3. Moves the "current selection" on the "ColourCode" sheet to F15
#     (without leaving that sheet tab-selected).
#  4. Leaves "AttendanceApproval" as the final active sheet/tab, with its
#     selection on G10 - this both sets sheetView tabSelected="1" there
#     and clears it from "AttendanceVerify" (previously the active tab),
#     and drives workbookView's activeTab down from 2 to 1.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new trailing worksheet -------------------------------------
$lastSheet = $wb.Sheets.Item($wb.Sheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet5"

# --- 2. Update the employee name + font colour on AttendanceApproval!B2 ----
$ws2 = $wb.Sheets.Item("AttendanceApproval")
$ws2.Range("B2").Value = "F Dinesh"
$ws2.Range("B2").Font.ColorIndex = 1

# --- 3. Re-point the ColourCode sheet's remembered selection to F15 --------
$ws4 = $wb.Sheets.Item("ColourCode")
$ws4.Range("F15").Select() | Out-Null

# --- 4. Finish with AttendanceApproval active/selected at G10 --------------
$ws2.Range("G10").Select() | Out-Null
